# Auto-generated cell updates replicating the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 37.5
$ws.Range("I6").Value = 37.5
$ws.Range("K6").Value = 112.5
$ws.Range("M6").Value = -0.5
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 8
$ws.Range("I8").Value = 5.6
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 16.8
$ws.Range("L8").Value = 60
$ws.Range("M8").Value = 122.2
$ws.Range("N8").Value = -338
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 1878.5
$ws.Range("I33").Value = 1810
$ws.Range("K33").Value = 1810
$ws.Range("M33").Value = -1581
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 4797.25
$ws.Range("J40").Value = 6844.5
$ws.Range("L40").Value = 6844.5
$ws.Range("N40").Value = -7194.5
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 5479.6816
$ws.Range("I41").Value = 721.4286
$ws.Range("J41").Value = 13806.625
$ws.Range("K41").Value = 721.4286
$ws.Range("L41").Value = 13806.625
$ws.Range("M41").Value = -281.4286
$ws.Range("N41").Value = -14686.625
# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 1402.125
$ws.Range("J88").Value = 513.8421
$ws.Range("L88").Value = 513.8421
$ws.Range("N88").Value = -1325.8421
# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 1402.125
$ws.Range("J91").Value = 513.8421
$ws.Range("L91").Value = 513.8421
$ws.Range("N91").Value = -3321.8421
# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 1841.2
$ws.Range("J112").Value = 1940.375
$ws.Range("L112").Value = 5821.125
$ws.Range("N112").Value = -8037.125

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 12893.186
$ws.Range("J2").Value = 1496.2727
$ws.Range("L2").Value = 1496.2727
$ws.Range("N2").Value = -1722.2727
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 3775.6316
$ws.Range("I61").Value = 1353.4546
$ws.Range("K61").Value = 1353.4546
$ws.Range("M61").Value = -1141.4546
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 12893.186
$ws.Range("J116").Value = 1496.2727
$ws.Range("L116").Value = 1496.2727
$ws.Range("N116").Value = -6084.2727
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 5495.25
$ws.Range("I122").Value = 3880.4443
$ws.Range("J122").Value = 7571.4287
$ws.Range("K122").Value = 11641.3329
$ws.Range("L122").Value = 22714.2861
$ws.Range("M122").Value = -9191.332900000001
$ws.Range("N122").Value = -27614.2861
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 3775.6316
$ws.Range("I136").Value = 1353.4546
$ws.Range("K136").Value = 4060.3638
$ws.Range("M136").Value = -1510.3638

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 12893.186
$ws.Range("J3").Value = 1496.2727
$ws.Range("L3").Value = 1496.2727
$ws.Range("N3").Value = -1724.2727
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 10072.5
$ws.Range("I22").Value = 10072.5
$ws.Range("K22").Value = 10072.5
$ws.Range("M22").Value = -9899.5
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2896.0386
$ws.Range("I86").Value = 1203.0714
$ws.Range("J86").Value = 4871.1665
$ws.Range("K86").Value = 1203.0714
$ws.Range("L86").Value = 4871.1665
$ws.Range("M86").Value = -80.07140000000004
$ws.Range("N86").Value = -7117.1665
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2896.0386
$ws.Range("I89").Value = 1203.0714
$ws.Range("J89").Value = 4871.1665
$ws.Range("K89").Value = 6015.357
$ws.Range("L89").Value = 24355.8325
$ws.Range("M89").Value = -399.357
$ws.Range("N89").Value = -35587.8325
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 32498.54
$ws.Range("I99").Value = 34790.25
$ws.Range("J99").Value = 4998
$ws.Range("K99").Value = 34790.25
$ws.Range("L99").Value = 4998
$ws.Range("M99").Value = -33292.25
$ws.Range("N99").Value = -7994
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2921.45
$ws.Range("I134").Value = 1465.5358
$ws.Range("K134").Value = 4396.607400000001
$ws.Range("M134").Value = -1861.607400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2677
$ws.Range("I31").Value = 1213.7142
$ws.Range("J31").Value = 5603.5713
$ws.Range("K31").Value = 1213.7142
$ws.Range("L31").Value = 5603.5713
$ws.Range("M31").Value = -918.7141999999999
$ws.Range("N31").Value = -6193.5713
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2677
$ws.Range("I34").Value = 1213.7142
$ws.Range("J34").Value = 5603.5713
$ws.Range("K34").Value = 1213.7142
$ws.Range("L34").Value = 5603.5713
$ws.Range("M34").Value = -1011.7142
$ws.Range("N34").Value = -6007.5713
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 409559.75
$ws.Range("I122").Value = 426333.12
$ws.Range("K122").Value = 1278999.36
$ws.Range("M122").Value = -1276549.36

$ws = $wb.Worksheets.Item("CUL")
# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 545.4545000000001
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -2909
$ws.Range("N46").Value = -1682
# Row 57 (Leve Item ID 4655)
$ws.Range("H57").Value = 136483
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 136483
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 409449
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -410567
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 10338.333
$ws.Range("I131").Value = 5343.3335
$ws.Range("J131").Value = 15333.333
$ws.Range("K131").Value = 16030.0005
$ws.Range("L131").Value = 45999.999
$ws.Range("M131").Value = -10990.0005
$ws.Range("N131").Value = -56079.999
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 41668028
$ws.Range("I140").Value = 41668028
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 125004084
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -124998904
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 64 (Leve Item ID 10640)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
# Row 67 (Leve Item ID 10640)
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 3211.7646
$ws.Range("I102").Value = 2206.8
$ws.Range("K102").Value = 2206.8
$ws.Range("M102").Value = -584.8000000000002
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 6449.5
$ws.Range("I122").Value = 6277.222
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 18831.666
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -16381.666
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 10628.286
$ws.Range("I40").Value = 12719
$ws.Range("J40").Value = 8727.637000000001
$ws.Range("K40").Value = 12719
$ws.Range("L40").Value = 8727.637000000001
$ws.Range("M40").Value = -12583
$ws.Range("N40").Value = -8999.637000000001
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 4409.294
$ws.Range("I122").Value = 3363.9092
$ws.Range("J122").Value = 6325.8335
$ws.Range("K122").Value = 10091.7276
$ws.Range("L122").Value = 18977.5005
$ws.Range("M122").Value = -7641.7276
$ws.Range("N122").Value = -23877.5005

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2377.3572
$ws.Range("I126").Value = 2080.2727
$ws.Range("K126").Value = 6240.8181
$ws.Range("M126").Value = -3770.8181

